$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.798.58'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.02%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.791.58'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.99%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '432.92'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +5.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.56'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +7.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.624'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.53%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.154'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -9.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000316'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -15.04%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.00'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +4.75%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.43'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +4.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.362.29'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.81%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.98'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.78%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.773.46'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.30%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.94'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.62%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.14'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +6.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '66.776.20'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.58%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '410.38'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.76'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.91%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +6.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.56'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '36.86'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.87%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.37'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +7.60%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.61'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +33.73%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.80'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +3.94%  '
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '13.85'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +11.21%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '713.99'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +5.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.134'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +10.34%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '41.79'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +8.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.13%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.152'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +26.23%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '56.13'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.85%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0476'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.86%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.78'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +43.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.92'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -5.78%  '
$ws.Range('B42').Value = 'PEPE'
$ws.Range('C42').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0₃0679'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -14.51%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.141'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.42%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.40%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.30'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +4.46%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.324'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +9.87%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.87%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.70'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +5.65%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '142.85'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -4.42%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.83'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.06%  '
